$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.019999999999999
$ws.Range("C2").Value = 1.037072668327167
$ws.Range("D2").Value = 1.043628341595469
$ws.Range("E2").Value = 1.045584069080439
$ws.Range("F2").Value = 1.055686291238313
$ws.Range("I2").Value = 1.036456979540641
$ws.Range("J2").Value = 1.042177528584906
$ws.Range("K2").Value = 1.046401671358399
$ws.Range("L2").Value = 1.048351899543528
$ws.Range("M2").Value = 1.058426066917074
$ws.Range("N2").Value = 1.04365754021493

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.038044225850713
$ws.Range("D3").Value = 1.044368873155201
$ws.Range("E3").Value = 1.046438048109888
$ws.Range("F3").Value = 1.056604036749747
$ws.Range("I3").Value = 1.036632779252859
$ws.Range("J3").Value = 1.042793202718466
$ws.Range("K3").Value = 1.046953460946603
$ws.Range("L3").Value = 1.049017237895863
$ws.Range("M3").Value = 1.059157037015363
$ws.Range("N3").Value = 1.044274088676378

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.038673192967832
$ws.Range("D4").Value = 1.044847929730826
$ws.Range("E4").Value = 1.046991237183872
$ws.Range("F4").Value = 1.057198397439465
$ws.Range("I4").Value = 1.036744728684612
$ws.Range("J4").Value = 1.043191281995167
$ws.Range("K4").Value = 1.047309722938608
$ws.Range("L4").Value = 1.049447720294752
$ws.Range("M4").Value = 1.059629909696378
$ws.Range("N4").Value = 1.044672733271313

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.038937682759817
$ws.Range("D5").Value = 1.045049295738069
$ws.Range("E5").Value = 1.047223942010007
$ws.Range("F5").Value = 1.057448389286197
$ws.Range("I5").Value = 1.036791359456251
$ws.Range("J5").Value = 1.043358561016328
$ws.Range("K5").Value = 1.047459306995139
$ws.Range("L5").Value = 1.049628685420865
$ws.Range("M5").Value = 1.059828676905069
$ws.Range("N5").Value = 1.044840249847871

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.038982095957593
$ws.Range("D6").Value = 1.045083104213219
$ws.Range("E6").Value = 1.047263022588626
$ws.Range("F6").Value = 1.057490371166817
$ws.Range("I6").Value = 1.036799163569412
$ws.Range("J6").Value = 1.04338664356177
$ws.Range("K6").Value = 1.047484411731843
$ws.Range("L6").Value = 1.04965906966379
$ws.Range("M6").Value = 1.05986204908948
$ws.Range("N6").Value = 1.044868372273748

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.038676726811489
$ws.Range("D7").Value = 1.044850620510475
$ws.Range("E7").Value = 1.046994346032427
$ws.Range("F7").Value = 1.057201737362729
$ws.Range("I7").Value = 1.036745353468655
$ws.Range("J7").Value = 1.043193517475672
$ws.Range("K7").Value = 1.047311722430089
$ws.Range("L7").Value = 1.049450138398969
$ws.Range("M7").Value = 1.059632565746932
$ws.Range("N7").Value = 1.044674971926457

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.037400947185703
$ws.Range("D8").Value = 1.043878631060245
$ws.Range("E8").Value = 1.04587254871274
$ws.Range("F8").Value = 1.055996339545038
$ws.Range("I8").Value = 1.03651676513278
$ws.Range("J8").Value = 1.042385660718052
$ws.Range("K8").Value = 1.046588312657826
$ws.Range("L8").Value = 1.048576760576839
$ws.Range("M8").Value = 1.05867312453158
$ws.Range("N8").Value = 1.043865967919578

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.035155230024456
$ws.Range("D9").Value = 1.042165027898826
$ws.Range("E9").Value = 1.043900508940992
$ws.Range("F9").Value = 1.053876302797093
$ws.Range("I9").Value = 1.036100172289784
$ws.Range("J9").Value = 1.040959834704433
$ws.Range("K9").Value = 1.045307624853312
$ws.Range("L9").Value = 1.047037527815955
$ws.Range("M9").Value = 1.056981647683585
$ws.Range("N9").Value = 1.042438117069476

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.033659727636023
$ws.Range("D10").Value = 1.041022146995464
$ws.Range("E10").Value = 1.042589059019807
$ws.Range("F10").Value = 1.052465734382275
$ws.Range("I10").Value = 1.03581321264756
$ws.Range("J10").Value = 1.040007806886366
$ws.Range("K10").Value = 1.044449897474266
$ws.Range("L10").Value = 1.046011279677142
$ws.Range("M10").Value = 1.055853516276417
$ws.Range("N10").Value = 1.041484737262708

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.033012557907566
$ws.Range("D11").Value = 1.040527170222608
$ws.Range("E11").Value = 1.042021971691876
$ws.Range("F11").Value = 1.051855620971629
$ws.Range("I11").Value = 1.035686774890858
$ws.Range("J11").Value = 1.039595229246128
$ws.Range("K11").Value = 1.044077571543895
$ws.Range("L11").Value = 1.045566893424691
$ws.Range("M11").Value = 1.055364923475025
$ws.Range("N11").Value = 1.041071573714897

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.032772229963359
$ws.Range("D12").Value = 1.040343299830167
$ws.Range("E12").Value = 1.041811448495112
$ws.Range("F12").Value = 1.051629099948681
$ws.Range("I12").Value = 1.035639483038133
$ws.Range("J12").Value = 1.039441928822678
$ws.Range("K12").Value = 1.043939135240614
$ws.Range("L12").Value = 1.045401827446493
$ws.Range("M12").Value = 1.055183423777465
$ws.Range("N12").Value = 1.040918055587261

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.032823778409825
$ws.Range("D13").Value = 1.040382741265104
$ws.Range("E13").Value = 1.041856601063051
$ws.Range("F13").Value = 1.051677684828606
$ws.Range("I13").Value = 1.035649642099291
$ws.Range("J13").Value = 1.039474814565193
$ws.Range("K13").Value = 1.043968836520838
$ws.Range("L13").Value = 1.045437234695073
$ws.Range("M13").Value = 1.055222356720431
$ws.Range("N13").Value = 1.040950988031303

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.03299269110257
$ws.Range("D14").Value = 1.040511971716224
$ws.Range("E14").Value = 1.042004567358657
$ws.Range("F14").Value = 1.051836894579572
$ws.Range("I14").Value = 1.035682872403523
$ws.Range("J14").Value = 1.039582558414962
$ws.Range("K14").Value = 1.044066131162762
$ws.Range("L14").Value = 1.045553249037045
$ws.Range("M14").Value = 1.055349920944316
$ws.Range("N14").Value = 1.041058884889698

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.033096771720656
$ws.Range("D15").Value = 1.040591593049
$ws.Range("E15").Value = 1.042095749987683
$ws.Range("F15").Value = 1.05193500255138
$ws.Range("I15").Value = 1.035703303347014
$ws.Range("J15").Value = 1.039648936261903
$ws.Range("K15").Value = 1.044126059362636
$ws.Range("L15").Value = 1.045624729180324
$ws.Range("M15").Value = 1.055428515574618
$ws.Range("N15").Value = 1.041125357000794

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.033702686493565
$ws.Range("D16").Value = 1.041054994941214
$ws.Range("E16").Value = 1.042626711258728
$ws.Range("F16").Value = 1.052506239869706
$ws.Range("I16").Value = 1.035821557979681
$ws.Range("J16").Value = 1.04003518114591
$ws.Range("K16").Value = 1.044474588130645
$ws.Range("L16").Value = 1.046040771916395
$ws.Range("M16").Value = 1.05588594046614
$ws.Range("N16").Value = 1.041512150396839

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.034082866456137
$ws.Range("D17").Value = 1.04134564838024
$ws.Range("E17").Value = 1.042959978703593
$ws.Range("F17").Value = 1.052864742785462
$ws.Range("I17").Value = 1.035895152159582
$ws.Range("J17").Value = 1.040277371051419
$ws.Range("K17").Value = 1.044692964255928
$ws.Range("L17").Value = 1.046301741461684
$ws.Range("M17").Value = 1.056172843839215
$ws.Range("N17").Value = 1.041754684239795

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.034304656835725
$ws.Range("D18").Value = 1.041515171703769
$ws.Range("E18").Value = 1.043154443059454
$ws.Range("F18").Value = 1.053073916224233
$ws.Range("I18").Value = 1.035937867793505
$ws.Range("J18").Value = 1.040418603143511
$ws.Range("K18").Value = 1.044820250105845
$ws.Range("L18").Value = 1.04645395916297
$ws.Range("M18").Value = 1.056340179569387
$ws.Range("N18").Value = 1.041896116897657

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.034380288020136
$ws.Range("D19").Value = 1.041572973055239
$ws.Range("E19").Value = 1.04322076307914
$ws.Range("F19").Value = 1.053145249905032
$ws.Range("I19").Value = 1.03595239697909
$ws.Range("J19").Value = 1.040466753982625
$ws.Range("K19").Value = 1.044863636122791
$ws.Range("L19").Value = 1.046505861195137
$ws.Range("M19").Value = 1.056397234972808
$ws.Range("N19").Value = 1.041944336116485

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.034042072824992
$ws.Range("D20").Value = 1.041314465032444
$ws.Range("E20").Value = 1.042924214484553
$ws.Range("F20").Value = 1.052826272108251
$ws.Range("I20").Value = 1.03588727797995
$ws.Range("J20").Value = 1.04025138979255
$ws.Range("K20").Value = 1.044669543788631
$ws.Range("L20").Value = 1.046273742039003
$ws.Range("M20").Value = 1.056142062864597
$ws.Range("N20").Value = 1.041728666084558

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.032942948870625
$ws.Range("D21").Value = 1.04047391694475
$ws.Range("E21").Value = 1.041960991697427
$ws.Range("F21").Value = 1.051790008442733
$ws.Range("I21").Value = 1.03567309593555
$ws.Range("J21").Value = 1.039550831929809
$ws.Range("K21").Value = 1.044037484113791
$ws.Range("L21").Value = 1.045519085728539
$ws.Range("M21").Value = 1.055312356853269
$ws.Range("N21").Value = 1.041027113349297

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.032252231823204
$ws.Range("D22").Value = 1.039945350762098
$ws.Range("E22").Value = 1.041356060275484
$ws.Range("F22").Value = 1.051139060650108
$ws.Range("I22").Value = 1.035536538182189
$ws.Range("J22").Value = 1.039110070382944
$ws.Range("K22").Value = 1.043639285896302
$ws.Range("L22").Value = 1.045044597249886
$ws.Range("M22").Value = 1.054790603967933
$ws.Range("N22").Value = 1.040585725870479

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.032618360982011
$ws.Range("D23").Value = 1.0402255608648
$ws.Range("E23").Value = 1.041676680530788
$ws.Range("F23").Value = 1.051484083814953
$ws.Range("I23").Value = 1.035609109255784
$ws.Range("J23").Value = 1.039343753905938
$ws.Range("K23").Value = 1.043850453541836
$ws.Range("L23").Value = 1.045296132825188
$ws.Range("M23").Value = 1.055067202614229
$ws.Range("N23").Value = 1.04081974125088

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.034060505594044
$ws.Range("D24").Value = 1.041328555475771
$ws.Range("E24").Value = 1.042940374563758
$ws.Range("F24").Value = 1.052843655151355
$ws.Range("I24").Value = 1.035890836633562
$ws.Range("J24").Value = 1.040263129707661
$ws.Range("K24").Value = 1.044680126765512
$ws.Range("L24").Value = 1.046286393777344
$ws.Range("M24").Value = 1.056155971493527
$ws.Range("N24").Value = 1.041740422671694

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.035735515862525
$ws.Range("D25").Value = 1.042608125735083
$ws.Range("E25").Value = 1.04440976241113
$ws.Range("F25").Value = 1.05442389737855
$ws.Range("I25").Value = 1.036209501567231
$ws.Range("J25").Value = 1.041328709295386
$ws.Range("K25").Value = 1.045639411613621
$ws.Range("L25").Value = 1.04743547724068
$ws.Range("M25").Value = 1.057419024762634
$ws.Range("N25").Value = 1.042807515504657
